$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
# Row 32
$ws.Cells.Item(32,8).Value = 448.5
$ws.Cells.Item(32,10).Value = 346.5
$ws.Cells.Item(32,12).Value = 346.5
$ws.Cells.Item(32,14).Value = -998.5

# Row 40
$ws.Cells.Item(40,8).Value = 2766.8667
$ws.Cells.Item(40,9).Value = 4171.4287
$ws.Cells.Item(40,10).Value = 1537.875
$ws.Cells.Item(40,11).Value = 4171.4287
$ws.Cells.Item(40,12).Value = 1537.875
$ws.Cells.Item(40,13).Value = -3996.4287
$ws.Cells.Item(40,14).Value = -1887.875

# Row 74
$ws.Cells.Item(74,8).Value = 4601.4287
$ws.Cells.Item(74,9).Value = 4601.4287
$ws.Cells.Item(74,11).Value = 4601.4287
$ws.Cells.Item(74,13).Value = -3665.4287

# Row 76
$ws.Cells.Item(76,8).Value = 3916.5
$ws.Cells.Item(76,9).Value = 3533.1333
$ws.Cells.Item(76,11).Value = 3533.1333
$ws.Cells.Item(76,13).Value = -3218.1333

# Row 77
$ws.Cells.Item(77,8).Value = 4601.4287
$ws.Cells.Item(77,9).Value = 4601.4287
$ws.Cells.Item(77,11).Value = 23007.1435
$ws.Cells.Item(77,13).Value = -18327.1435

# Row 79
$ws.Cells.Item(79,8).Value = 3916.5
$ws.Cells.Item(79,9).Value = 3533.1333
$ws.Cells.Item(79,11).Value = 3533.1333
$ws.Cells.Item(79,13).Value = -2441.1333

# Row 100
$ws.Cells.Item(100,8).Value = 1563
$ws.Cells.Item(100,9).Value = 1781.1111
$ws.Cells.Item(100,10).Value = 1344.8889
$ws.Cells.Item(100,11).Value = 1781.1111
$ws.Cells.Item(100,12).Value = 1344.8889
$ws.Cells.Item(100,13).Value = -1240.1111
$ws.Cells.Item(100,14).Value = -2426.8889

# Row 129
$ws.Cells.Item(129,8).Value = 882.5106
$ws.Cells.Item(129,10).Value = 922.36365
$ws.Cells.Item(129,12).Value = 2767.09095
$ws.Cells.Item(129,14).Value = -12767.09095

# Row 132
$ws.Cells.Item(132,8).Value = 777.2414
$ws.Cells.Item(132,9).Value = 777.2414
$ws.Cells.Item(132,10).Value = 0
$ws.Cells.Item(132,11).Value = 2331.7242
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,13).Value = 198.2757999999999
$ws.Cells.Item(132,14).Value = $null

# Row 137
$ws.Cells.Item(137,8).Value = 643739.8
$ws.Cells.Item(137,9).Value = 3497.5881
$ws.Cells.Item(137,10).Value = 954714.5600000001
$ws.Cells.Item(137,11).Value = 10492.7643
$ws.Cells.Item(137,12).Value = 2864143.68
$ws.Cells.Item(137,13).Value = -7942.764299999999
$ws.Cells.Item(137,14).Value = -2869243.68

# Row 138
$ws.Cells.Item(138,8).Value = 3716.2031
$ws.Cells.Item(138,9).Value = 2518.4546
$ws.Cells.Item(138,11).Value = 7555.3638
$ws.Cells.Item(138,13).Value = -2415.3638

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32,8).Value = 17196.447
$ws.Cells.Item(32,9).Value = 18018.129
$ws.Cells.Item(32,11).Value = 18018.129
$ws.Cells.Item(32,13).Value = -17731.129

# Row 97
$ws.Cells.Item(97,8).Value = 913.5172
$ws.Cells.Item(97,9).Value = 867.3182
$ws.Cells.Item(97,10).Value = 1058.7142
$ws.Cells.Item(97,11).Value = 867.3182
$ws.Cells.Item(97,12).Value = 1058.7142
$ws.Cells.Item(97,13).Value = -371.3182
$ws.Cells.Item(97,14).Value = -2050.7142

# Row 102
$ws.Cells.Item(102,8).Value = 3675.8333
$ws.Cells.Item(102,9).Value = 3434.4443
$ws.Cells.Item(102,11).Value = 3434.4443
$ws.Cells.Item(102,13).Value = -1812.4443

# Row 132
$ws.Cells.Item(132,8).Value = 3162.7932
$ws.Cells.Item(132,9).Value = 2884.2632
$ws.Cells.Item(132,10).Value = 3692
$ws.Cells.Item(132,11).Value = 8652.7896
$ws.Cells.Item(132,12).Value = 11076
$ws.Cells.Item(132,13).Value = -6122.7896
$ws.Cells.Item(132,14).Value = -16136

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
# Row 94
$ws.Cells.Item(94,8).Value = 1584.2609
$ws.Cells.Item(94,9).Value = 1560.4706
$ws.Cells.Item(94,10).Value = 1651.6666
$ws.Cells.Item(94,11).Value = 1560.4706
$ws.Cells.Item(94,12).Value = 1651.6666
$ws.Cells.Item(94,13).Value = -1109.4706
$ws.Cells.Item(94,14).Value = -2553.6666

# Row 99
$ws.Cells.Item(99,8).Value = 3158.9
$ws.Cells.Item(99,9).Value = 1227
$ws.Cells.Item(99,10).Value = 7666.6665
$ws.Cells.Item(99,11).Value = 1227
$ws.Cells.Item(99,12).Value = 7666.6665
$ws.Cells.Item(99,13).Value = 271
$ws.Cells.Item(99,14).Value = -10662.6665

# Row 105
$ws.Cells.Item(105,8).Value = 5800.2
$ws.Cells.Item(105,9).Value = 6334
$ws.Cells.Item(105,10).Value = 4999.5
$ws.Cells.Item(105,11).Value = 6334
$ws.Cells.Item(105,12).Value = 4999.5
$ws.Cells.Item(105,13).Value = -4587
$ws.Cells.Item(105,14).Value = -8493.5

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31,8).Value = 1097760.4
$ws.Cells.Item(31,9).Value = 18937.545
$ws.Cells.Item(31,10).Value = 1613719
$ws.Cells.Item(31,11).Value = 18937.545
$ws.Cells.Item(31,12).Value = 1613719
$ws.Cells.Item(31,13).Value = -18642.545
$ws.Cells.Item(31,14).Value = -1614309

# Row 34
$ws.Cells.Item(34,8).Value = 1097760.4
$ws.Cells.Item(34,9).Value = 18937.545
$ws.Cells.Item(34,10).Value = 1613719
$ws.Cells.Item(34,11).Value = 18937.545
$ws.Cells.Item(34,12).Value = 1613719
$ws.Cells.Item(34,13).Value = -18735.545
$ws.Cells.Item(34,14).Value = -1614123

# Row 58
$ws.Cells.Item(58,8).Value = 1937873.1
$ws.Cells.Item(58,9).Value = 3368232
$ws.Cells.Item(58,10).Value = 6888.5
$ws.Cells.Item(58,11).Value = 3368232
$ws.Cells.Item(58,12).Value = 6888.5
$ws.Cells.Item(58,13).Value = -3368029
$ws.Cells.Item(58,14).Value = -7294.5

# Row 62
$ws.Cells.Item(62,8).Value = 3699.9
$ws.Cells.Item(62,9).Value = 3665.8333
$ws.Cells.Item(62,11).Value = 3665.8333
$ws.Cells.Item(62,13).Value = -3041.8333

# Row 65
$ws.Cells.Item(65,8).Value = 3699.9
$ws.Cells.Item(65,9).Value = 3665.8333
$ws.Cells.Item(65,11).Value = 18329.1665
$ws.Cells.Item(65,13).Value = -15209.1665

# Row 132
$ws.Cells.Item(132,8).Value = 2743.5122
$ws.Cells.Item(132,9).Value = 2499.1516
$ws.Cells.Item(132,10).Value = 3751.5
$ws.Cells.Item(132,11).Value = 7497.4548
$ws.Cells.Item(132,12).Value = 11254.5
$ws.Cells.Item(132,13).Value = -4967.4548
$ws.Cells.Item(132,14).Value = -16314.5

# Row 134
$ws.Cells.Item(134,8).Value = 2413.75
$ws.Cells.Item(134,9).Value = 1926.75
$ws.Cells.Item(134,10).Value = 3387.75
$ws.Cells.Item(134,11).Value = 5780.25
$ws.Cells.Item(134,12).Value = 10163.25
$ws.Cells.Item(134,13).Value = -3245.25
$ws.Cells.Item(134,14).Value = -15233.25

# Row 136
$ws.Cells.Item(136,8).Value = 1937873.1
$ws.Cells.Item(136,9).Value = 3368232
$ws.Cells.Item(136,10).Value = 6888.5
$ws.Cells.Item(136,11).Value = 10104696
$ws.Cells.Item(136,12).Value = 20665.5
$ws.Cells.Item(136,13).Value = -10102146
$ws.Cells.Item(136,14).Value = -25765.5

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
# Row 107
$ws.Cells.Item(107,8).Value = 568.6869
$ws.Cells.Item(107,10).Value = 678.451
$ws.Cells.Item(107,12).Value = 2035.353
$ws.Cells.Item(107,14).Value = -5875.353

# Row 131
$ws.Cells.Item(131,8).Value = 1266.878
$ws.Cells.Item(131,10).Value = 1086.1936
$ws.Cells.Item(131,12).Value = 3258.5808
$ws.Cells.Item(131,14).Value = -13338.5808

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Cells.Item(70,8).Value = 5700.1396
$ws.Cells.Item(70,9).Value = 5477.154
$ws.Cells.Item(70,11).Value = 5477.154
$ws.Cells.Item(70,13).Value = -5207.154

# Row 73
$ws.Cells.Item(73,8).Value = 5700.1396
$ws.Cells.Item(73,9).Value = 5477.154
$ws.Cells.Item(73,11).Value = 5477.154
$ws.Cells.Item(73,13).Value = -4541.154

# Row 97
$ws.Cells.Item(97,8).Value = 1320.8235
$ws.Cells.Item(97,9).Value = 1470.25
$ws.Cells.Item(97,10).Value = 962.2
$ws.Cells.Item(97,11).Value = 1470.25
$ws.Cells.Item(97,12).Value = 962.2
$ws.Cells.Item(97,13).Value = -974.25
$ws.Cells.Item(97,14).Value = -1954.2

# Row 102
$ws.Cells.Item(102,8).Value = 4211.2383
$ws.Cells.Item(102,9).Value = 3771.2
$ws.Cells.Item(102,11).Value = 3771.2
$ws.Cells.Item(102,13).Value = -2149.2

# Row 122
$ws.Cells.Item(122,8).Value = 10901.5
$ws.Cells.Item(122,9).Value = 15687.143
$ws.Cells.Item(122,11).Value = 47061.429
$ws.Cells.Item(122,13).Value = -44611.429

# Row 132
$ws.Cells.Item(132,8).Value = 13423.053
$ws.Cells.Item(132,9).Value = 13602.667
$ws.Cells.Item(132,11).Value = 40808.001
$ws.Cells.Item(132,13).Value = -38278.001

# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
# Row 93
$ws.Cells.Item(93,8).Value = 834.8823
$ws.Cells.Item(93,9).Value = 837.0625
$ws.Cells.Item(93,11).Value = 837.0625
$ws.Cells.Item(93,13).Value = 410.9375

# Row 122
$ws.Cells.Item(122,8).Value = 5997.625
$ws.Cells.Item(122,9).Value = 5811.7646
$ws.Cells.Item(122,10).Value = 7050.8335
$ws.Cells.Item(122,11).Value = 17435.2938
$ws.Cells.Item(122,12).Value = 21152.5005
$ws.Cells.Item(122,13).Value = -14985.2938
$ws.Cells.Item(122,14).Value = -26052.5005

# Row 136
$ws.Cells.Item(136,8).Value = 4110.483
$ws.Cells.Item(136,9).Value = 3036.8462
$ws.Cells.Item(136,10).Value = 4982.8125
$ws.Cells.Item(136,11).Value = 9110.5386
$ws.Cells.Item(136,12).Value = 14948.4375
$ws.Cells.Item(136,13).Value = -6560.5386
$ws.Cells.Item(136,14).Value = -20048.4375

# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
# Row 96
$ws.Cells.Item(96,8).Value = 1429.3334
$ws.Cells.Item(96,9).Value = 1498
$ws.Cells.Item(96,10).Value = 1395
$ws.Cells.Item(96,11).Value = 1498
$ws.Cells.Item(96,12).Value = 1395
$ws.Cells.Item(96,13).Value = -125
$ws.Cells.Item(96,14).Value = -4141

# Row 122
$ws.Cells.Item(122,8).Value = 2500
